$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of column R into column S for rows 3-14 so the new
# column inherits the same number formats / styles as the existing data.
$ws.Range("R3:R14").Copy() | Out-Null
$ws.Range("S3:S14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new "2023" column of data.
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 75.3
$ws.Range("S5").Value = 76.8
$ws.Range("S6").Value = 1440
$ws.Range("S7").Value = 1387
$ws.Range("S8").Value = 1219.7
$ws.Range("S9").Value = 37.5
$ws.Range("S10").Value = 20.4
$ws.Range("S11").Value = 39.3
$ws.Range("S12").Value = 0.1
$ws.Range("S13").Value = 37.6
$ws.Range("S14").Value = 0.3841

# Update the active selection, matching the saved sheet view state.
$ws.Range("I22").Select() | Out-Null
